$wb = $excel.ActiveWorkbook

# --- "About" sheet updates -------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Update currency-year explanation text: 2015 -> 2012
$about.Range("A26").Value = "We adjust the sources' dollars to 2012 dollars using the following conversion factors:"

# Update conversion-factor note text: 2015 -> 2012, and the USD->CAD adjustment factor is gone
$about.Range("A27").Value = 1.278
$about.Range("B27").Value = "2002 to 2012, for U.S. Forest Service (2006) ""Regional Cost Information…"""

# Row 28 used to hold the 2015 USD->CAD conversion factor (1.3901) and its unit
# label; that conversion no longer exists. Instead, row 28 now carries the
# "see cpi.xlsx" note that used to live in row 30, and row 30 is removed.
$about.Range("A28").Value = $about.Range("A30").Value()
$about.Range("B28").ClearContents()
$about.Range("A30").ClearContents()

# --- "Forest Mgmt Costs" sheet updates --------------------------------------
$fmc = $wb.Worksheets.Item("Forest Mgmt Costs")

# Unit labels simplified (no more USD/CAD distinction since conversion removed)
$fmc.Range("B39").Value = "2002$ / acre"
$fmc.Range("B40").Value = "2012$ / acre"

# Drop the USD->CAD multiplier (About!A28) from the final cost formula
$fmc.Range("A40").Formula = "=A39*About!A27"

# The "About" sheet becomes the active/selected tab again
$about.Activate()

$wb.Application.CalculateFull()
